# Update Name of Algo
# Apply updated values to result_data_KNN.xlsx per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.183
$ws.Range("D5").Value = -8.148999999999999
$ws.Range("D9").Value = -7.746
$ws.Range("D11").Value = -8.359
$ws.Range("A21").Value = -20.959
$ws.Range("D21").Value = -7.897
$ws.Range("A23").Value = -21.368
$ws.Range("A25").Value = -21.937
